# Add two new columns, I ("I0") and J ("IF"), to the data table.
# Headers go in row 1 (same style as the existing headers, e.g. H1),
# data values fill rows 2-50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells -----------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting already used by the other header cells (bold,
# centered, thin border) by copying the style from an existing header.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Data values --------------------------------------------------------
$iValues = @(6,7,8,6,5,10,9,4,6,6,7,7,6,8,6,7,7,9,4,6,8,9,6,9,7,7,9,9,9,8,8,4,7,8,6,5,6,8,6,9,6,6,8,6,6,4,7,7,6)
$jValues = @(6,8,8,6,6,10,9,5,6,7,8,8,7,8,6,8,8,9,5,7,8,9,7,9,7,7,9,9,9,9,8,5,8,9,6,6,6,8,7,9,7,6,8,6,6,4,7,7,6)

for ($r = 2; $r -le 50; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]   # column I
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]  # column J
}
